$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.072.79"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  -1.46%  "
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').Value = "'2.940.47"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  -2.35%  "
$ws.Range('E3').Style = "Normal"

$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = "'  -0.09%  "
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').Value = "'376.02"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  -1.47%  "
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').Value = "'102.43"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'  -4.14%  "
$ws.Range('E6').Style = "Normal"

$ws.Range('D7').Value = "'0.536"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = "'  -2.34%  "
$ws.Range('E7').Style = "Normal"

$ws.Range('E8').Value = "'  -0.01%  "
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').Value = "'0.584"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  -3.23%  "
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').Value = "'36.67"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "'  -3.06%  "
$ws.Range('E10').Style = "Normal"

$ws.Range('E11').Value = "'  -1.12%  "
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').Value = "'0.0835"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "'  -1.54%  "
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').Value = "'3.397.25"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'  -2.62%  "
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').Value = "'17.95"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  -4.65%  "
$ws.Range('E14').Style = "Normal"

$ws.Range('E15').Value = "'  -3.21%  "
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').Value = "'2.927.17"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  -2.96%  "
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').Value = "'0.973"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  -0.26%  "
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').Value = "'50.946.52"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "'  -1.80%  "
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').Value = "'3.17"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "'  -6.89%  "
$ws.Range('E19').Style = "Normal"

$ws.Range('E20').Value = "'  -4.72%  "
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').Value = "'12.48"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  -5.09%  "
$ws.Range('E21').Style = "Normal"

$ws.Range('D22').Value = "'0.0₃0948"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = "'  -1.67%  "
$ws.Range('E22').Style = "Normal"

$ws.Range('B23').Value = "'BitcoinCash"
$ws.Range('B23').Style = "Normal"
$ws.Range('C23').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('C23').Style = "Normal"
$ws.Range('D23').Value = "'263.00"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "'  -0.74%  "
$ws.Range('E23').Style = "Normal"

$ws.Range('B24').Value = "'Litecoin"
$ws.Range('B24').Style = "Normal"
$ws.Range('C24').Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range('C24').Style = "Normal"
$ws.Range('D24').Value = "'68.18"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "'  -1.29%  "
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').Value = "'2.85"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  +1.49%  "
$ws.Range('E25').Style = "Normal"

$ws.Range('D26').Value = "'8.17"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "'  +8.70%  "
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').Value = "'7.73"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'  +6.45%  "
$ws.Range('E27').Style = "Normal"

$ws.Range('E28').Value = "'  -3.44%  "
$ws.Range('E28').Style = "Normal"

$ws.Range('E29').Value = "'  -0.01%  "
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').Value = "'0.112"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "'  +3.85%  "
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').Value = "'25.65"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'  -2.30%  "
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').Value = "'9.84"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'  -1.60%  "
$ws.Range('E32').Style = "Normal"

$ws.Range('B33').Value = "'InjectiveProtocol"
$ws.Range('B33').Style = "Normal"
$ws.Range('C33').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C33').Style = "Normal"
$ws.Range('D33').Value = "'34.06"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'  -2.94%  "
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').Value = "'50.65"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  -1.33%  "
$ws.Range('E34').Style = "Normal"

$ws.Range('B35').Value = "'VeChain"
$ws.Range('B35').Style = "Normal"
$ws.Range('C35').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C35').Style = "Normal"
$ws.Range('D35').Value = "'0.0454"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'  -0.95%  "
$ws.Range('E35').Style = "Normal"

$ws.Range('E36').Value = "'  -2.99%  "
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').Value = "'  -0.24%  "
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').Value = "'2.98"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'  -5.55%  "
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').Value = "'2.58"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'  -2.69%  "
$ws.Range('E39').Style = "Normal"

$ws.Range('B40').Value = "'Celestia"
$ws.Range('B40').Style = "Normal"
$ws.Range('C40').Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range('C40').Style = "Normal"
$ws.Range('D40').Value = "'16.45"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  -6.97%  "
$ws.Range('E40').Style = "Normal"

$ws.Range('B41').Value = "'Stellar"
$ws.Range('B41').Style = "Normal"
$ws.Range('C41').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('C41').Style = "Normal"
$ws.Range('D41').Value = "'0.114"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = "'  -2.19%  "
$ws.Range('E41').Style = "Normal"

$ws.Range('E42').Value = "'  -5.24%  "
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').Value = "'121.24"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "'  -2.73%  "
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').Value = "'21.16"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'  -6.31%  "
$ws.Range('E44').Style = "Normal"

$ws.Range('D45').Value = "'2.05"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'  -1.62%  "
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').Value = "'  -2.81%  "
$ws.Range('E46').Style = "Normal"

$ws.Range('E47').Value = "'  -1.61%  "
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').Value = "'2.007.02"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "'  -2.72%  "
$ws.Range('E48').Style = "Normal"

$ws.Range('D49').Value = "'3.23"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "'  -2.50%  "
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').Value = "'0.0347"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  -2.81%  "
$ws.Range('E50').Style = "Normal"

$ws.Range('E51').Value = "'  -4.23%  "
$ws.Range('E51').Style = "Normal"
